# New command to delete action file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "V2 Command" sheet (rId1)

# --- Hide helper columns D (Index) and E (Full record as Header) ---
$ws.Range("D1:E1").EntireColumn.Hidden = $true

# --- Populate row 42 with the new "Delete action file" command ---
$ws.Range("B42").Value = 75
$ws.Range("C42").Value = "Delete action file"
$ws.Range("D42").Value = "Yes {3}"
$ws.Range("E42").Value = "Action ID"
$ws.Range("F42").Value = "A9 9A 03 75 00 78 ED"
$ws.Range("F42").WrapText = $false
$ws.Range("I42").Value = "Yes {7}"
$ws.Range("J42").Value = "Yes {1}"
$ws.Range("K42").Value = "{result}"

# --- Update view state: scroll/selection on the sheet ---
$ws.Activate()
$win = $excel.ActiveWindow
$ws.Range("B39:F42").Select()
$win.ScrollRow = 15
$win.ScrollColumn = 4
